$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2359.8572
$ws.Range("I53").Value = 2331.9285
$ws.Range("J53").Value = 2387.7856
$ws.Range("K53").Value = 2331.9285
$ws.Range("L53").Value = 2387.7856
$ws.Range("M53").Value = -1694.9285
$ws.Range("N53").Value = -3661.7856
$ws.Range("H62").Value = 100029880
$ws.Range("I62").Value = 142861740
$ws.Range("K62").Value = 142861740
$ws.Range("M62").Value = -142861116
$ws.Range("H65").Value = 100029880
$ws.Range("I65").Value = 142861740
$ws.Range("K65").Value = 714308700
$ws.Range("M65").Value = -714305580
$ws.Range("H132").Value = 1876.8889
$ws.Range("I132").Value = 1557.258
$ws.Range("J132").Value = 3858.6
$ws.Range("K132").Value = 4671.774
$ws.Range("L132").Value = 11575.8
$ws.Range("M132").Value = -2141.774
$ws.Range("N132").Value = -16635.8
$ws.Range("H138").Value = 3685.5103
$ws.Range("I138").Value = 907.5714
$ws.Range("J138").Value = 5768.9644
$ws.Range("K138").Value = 2722.7142
$ws.Range("L138").Value = 17306.8932
$ws.Range("M138").Value = 2417.2858
$ws.Range("N138").Value = -27586.8932
$ws.Range("H141").Value = 2695
$ws.Range("I141").Value = 2990
$ws.Range("K141").Value = 8970
$ws.Range("M141").Value = -3790

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1760.25
$ws.Range("I32").Value = 1666.5676
$ws.Range("J32").Value = 5226.5
$ws.Range("K32").Value = 1666.5676
$ws.Range("L32").Value = 5226.5
$ws.Range("M32").Value = -1379.5676
$ws.Range("N32").Value = -5800.5
$ws.Range("H61").Value = 4965.88
$ws.Range("I61").Value = 2298.7646
$ws.Range("J61").Value = 10633.5
$ws.Range("K61").Value = 2298.7646
$ws.Range("L61").Value = 10633.5
$ws.Range("M61").Value = -2086.7646
$ws.Range("N61").Value = -11057.5
$ws.Range("H80").Value = 30000
$ws.Range("I80").Value = 30000
$ws.Range("K80").Value = 30000
$ws.Range("M80").Value = -29002
$ws.Range("H83").Value = 30000
$ws.Range("I83").Value = 30000
$ws.Range("K83").Value = 90000
$ws.Range("M83").Value = -85008
$ws.Range("H102").Value = 1341.3334
$ws.Range("I102").Value = 1007.4
$ws.Range("K102").Value = 1007.4
$ws.Range("M102").Value = 614.6
$ws.Range("H109").Value = 59378
$ws.Range("J109").Value = 59378
$ws.Range("L109").Value = 59378
$ws.Range("N109").Value = -62152
$ws.Range("H132").Value = 9086.929
$ws.Range("I132").Value = 7945.6
$ws.Range("K132").Value = 23836.8
$ws.Range("M132").Value = -21306.8
$ws.Range("H136").Value = 4965.88
$ws.Range("I136").Value = 2298.7646
$ws.Range("J136").Value = 10633.5
$ws.Range("K136").Value = 6896.293799999999
$ws.Range("L136").Value = 31900.5
$ws.Range("M136").Value = -4346.293799999999
$ws.Range("N136").Value = -37000.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 18520216
$ws.Range("I20").Value = 20835092
$ws.Range("J20").Value = 1209
$ws.Range("K20").Value = 20835092
$ws.Range("L20").Value = 1209
$ws.Range("M20").Value = -20834845
$ws.Range("N20").Value = -1703
$ws.Range("H105").Value = 3291.125
$ws.Range("I105").Value = 2899.6316
$ws.Range("J105").Value = 4778.8
$ws.Range("K105").Value = 2899.6316
$ws.Range("L105").Value = 4778.8
$ws.Range("M105").Value = -1152.6316
$ws.Range("N105").Value = -8272.799999999999
$ws.Range("H128").Value = 2690.7144
$ws.Range("I128").Value = 2690.7144
$ws.Range("K128").Value = 8072.1432
$ws.Range("M128").Value = -5582.1432

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10023.579
$ws.Range("I31").Value = 4030.5833
$ws.Range("K31").Value = 4030.5833
$ws.Range("M31").Value = -3735.5833
$ws.Range("H34").Value = 10023.579
$ws.Range("I34").Value = 4030.5833
$ws.Range("K34").Value = 4030.5833
$ws.Range("M34").Value = -3828.5833
$ws.Range("H132").Value = 6876.375
$ws.Range("I132").Value = 2342.5
$ws.Range("K132").Value = 7027.5
$ws.Range("M132").Value = -4497.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1400
$ws.Range("I3").Value = 1400
$ws.Range("K3").Value = 4200
$ws.Range("M3").Value = -4088
$ws.Range("H131").Value = 1889.3334
$ws.Range("J131").Value = 2164.24
$ws.Range("L131").Value = 6492.719999999999
$ws.Range("N131").Value = -16572.72

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 11823.333
$ws.Range("I99").Value = 12735
$ws.Range("K99").Value = 12735
$ws.Range("M99").Value = -10489
$ws.Range("H113").Value = 5442.049
$ws.Range("I113").Value = 2029.4445
$ws.Range("J113").Value = 8112.7827
$ws.Range("K113").Value = 2029.4445
$ws.Range("L113").Value = 8112.7827
$ws.Range("M113").Value = 140.5554999999999
$ws.Range("N113").Value = -12452.7827
$ws.Range("H126").Value = 3806
$ws.Range("I126").Value = 3806
$ws.Range("K126").Value = 11418
$ws.Range("M126").Value = -8948

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2431.6924
$ws.Range("J22").Value = 2384.3333
$ws.Range("L22").Value = 2384.3333
$ws.Range("N22").Value = -2974.3333
$ws.Range("H27").Value = 2431.6924
$ws.Range("J27").Value = 2384.3333
$ws.Range("L27").Value = 2384.3333
$ws.Range("N27").Value = -2598.3333
$ws.Range("H61").Value = 8050
$ws.Range("J61").Value = 8050
$ws.Range("L61").Value = 8050
$ws.Range("N61").Value = -8454
$ws.Range("H100").Value = 3724.8823
$ws.Range("I100").Value = 3057.8
$ws.Range("K100").Value = 3057.8
$ws.Range("M100").Value = -2516.8
$ws.Range("H113").Value = 8050
$ws.Range("J113").Value = 8050
$ws.Range("L113").Value = 8050
$ws.Range("N113").Value = -12390
$ws.Range("H122").Value = 3570.9788
$ws.Range("I122").Value = 2865.0278
$ws.Range("K122").Value = 8595.0834
$ws.Range("M122").Value = -6145.0834
$ws.Range("H132").Value = 25013300
$ws.Range("J132").Value = 13000
$ws.Range("L132").Value = 39000
$ws.Range("N132").Value = -44060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22231130
$ws.Range("I81").Value = 2881
$ws.Range("J81").Value = 100030000
$ws.Range("K81").Value = 5762
$ws.Range("L81").Value = 200060000
$ws.Range("M81").Value = -4701
$ws.Range("N81").Value = -200062122
$ws.Range("H84").Value = 22231130
$ws.Range("I84").Value = 2881
$ws.Range("J84").Value = 100030000
$ws.Range("K84").Value = 28810
$ws.Range("L84").Value = 1000300000
$ws.Range("M84").Value = -23506
$ws.Range("N84").Value = -1000310608
$ws.Range("H107").Value = 23810336
$ws.Range("I107").Value = 866.5
$ws.Range("J107").Value = 41667436
$ws.Range("K107").Value = 2599.5
$ws.Range("L107").Value = 125002308
$ws.Range("M107").Value = -679.5
$ws.Range("N107").Value = -125006148
$ws.Range("H122").Value = 151304.44
$ws.Range("I122").Value = 236660
$ws.Range("J122").Value = 6200
$ws.Range("K122").Value = 709980
$ws.Range("L122").Value = 18600
$ws.Range("M122").Value = -707530
$ws.Range("N122").Value = -23500
$ws.Range("H132").Value = 14711883
$ws.Range("I132").Value = 23815926
$ws.Range("J132").Value = 5350.846
$ws.Range("K132").Value = 71447778
$ws.Range("L132").Value = 16052.538
$ws.Range("M132").Value = -71445248
$ws.Range("N132").Value = -21112.538
$ws.Range("H133").Value = 128799.8
$ws.Range("J133").Value = 128799.8
$ws.Range("L133").Value = 128799.8
$ws.Range("N133").Value = -138919.8
$ws.Range("H136").Value = 25028990
$ws.Range("I136").Value = 66667364
$ws.Range("K136").Value = 200002092
$ws.Range("M136").Value = -199999542
